# Update the "修改时间" (last modified time) timestamp column on each
# portfolio sheet from 202509211526 to 202509211530.
#
# Sheet 1 "大智投资组合"   -> column E, data rows 2-9
# Sheet 2 "大成投资组合"   -> column E, data rows 2-11
# Sheet 3 "我的投资组合"   -> column G, data rows 2-13
#
# The timestamp values are stored as text (not numbers), so we force the
# target cells to text format before writing the new value to avoid Excel
# auto-converting the numeric-looking string into a real number.

$wb = $excel.ActiveWorkbook

$oldTimestamp = "202509211526"
$newTimestamp = "202509211530"

$sheetConfig = @(
    @{ Name = "大智投资组合"; Column = "E"; LastRow = 9 },
    @{ Name = "大成投资组合"; Column = "E"; LastRow = 11 },
    @{ Name = "我的投资组合"; Column = "G"; LastRow = 13 }
)

foreach ($cfg in $sheetConfig) {
    $ws = $wb.Worksheets.Item($cfg.Name)
    for ($row = 2; $row -le $cfg.LastRow; $row++) {
        $cell = $ws.Range($cfg.Column + $row)
        if ($cell.Value2 -eq $oldTimestamp) {
            $cell.NumberFormat = "@"
            $cell.Value = $newTimestamp
        }
    }
}
